$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "$ 38 , 656 , 710 , 449 , 061.55 AND 26 %"
$ws.Range("E2").Value = "17 YEARS AND 16.30 %"
$ws.Range("H2").Value = "$ 28 , 605 , 965 , 732 , 305.55 AND $ 131 , 685 , 600 , 094.16"
$ws.Range("I2").Value = "$ 67 , 186 , 820 , 273 , 079.76 AND $ 42 , 595 , 191 , 764.52"
